$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - the table layout is changing (columns reordered, new audit
# columns added, and row content refreshed), so clear existing content
# first and rewrite the whole table.
$ws.Cells.ClearContents()

# Header row
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "descr"
$ws.Range("D1").Value = "lang_code"
$ws.Range("E1").Value = "is_active"
$ws.Range("F1").Value = "cr_by"
$ws.Range("G1").Value = "cr_dtimes"
$ws.Range("H1").Value = "upd_by"
$ws.Range("I1").Value = "upd_dtimes"
$ws.Range("J1").Value = "is_deleted"
$ws.Range("K1").Value = "del_dtimes"

$rows = @(
    @{ id = 10001; name = "PrÃ©-Enregistrement"; descr = "Portail Web pour les prÃ©-enregistrements" },
    @{ id = 10002; name = "Enregistrement"; descr = "Application pour les enregistrements" },
    @{ id = 10003; name = "Traitement"; descr = "Application pour les traitements post-enregistrements" },
    @{ id = 10004; name = "Authentification"; descr = "Application pour l'authentification des fournisseurs de services" },
    @{ id = 10005; name = "Administration"; descr = "Portail Web pour la configuration des applications" },
    @{ id = 10006; name = "Portail RÃ©sident"; descr = "Portail Web pour les services dÃ©diÃ©s aux rÃ©sidents" },
    @{ id = 10007; name = "Processeur dinscription"; descr = "Demande de post-inscription" },
    @{ id = 10008; name = "Commune"; descr = "Commune" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.id
    $ws.Cells.Item($r, 2).Value = $row.name
    $ws.Cells.Item($r, 3).Value = $row.descr
    $ws.Cells.Item($r, 4).Value = "fra"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = 45079.577192083336
    $ws.Cells.Item($r, 7).NumberFormat = "mm:ss.0"
    $ws.Cells.Item($r, 8).Value = "NULL"
    $ws.Cells.Item($r, 9).Value = "NULL"
    $ws.Cells.Item($r, 10).Value = $false
    $ws.Cells.Item($r, 11).Value = "NULL"
    $r++
}

$ws.Range("D11").Select() | Out-Null
